# Updates cryptos list: refreshed prices/volume percentages, and restores
# the original (swapped) row order for Toncoin/Cosmos, EthereumClassic/Monero,
# and Algorand/InjectiveProtocol.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.917.06'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').Value = '2.350.63'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''241.05'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('D6').Value = '''0.673'
$ws.Range('E6').Value = '  -3.18%  '
$ws.Range('D7').Value = '''72.56'
$ws.Range('E7').Value = '  -5.32%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '''0.598'
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('E10').Value = '  -3.09%  '
$ws.Range('D11').Value = '''58.36'
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').Value = '''33.12'
$ws.Range('E12').Value = '  +2.58%  '
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').Value = '''7.27'
$ws.Range('E14').Value = '  -3.28%  '
$ws.Range('D15').Value = '2.700.42'
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('D16').Value = '''16.34'
$ws.Range('E16').Value = '  -5.14%  '
$ws.Range('E17').Value = '  -2.59%  '
$ws.Range('D18').Value = '2.351.64'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').Value = '43.831.80'
$ws.Range('E19').Value = '  -1.58%  '
$ws.Range('E20').Value = '  -1.04%  '
$ws.Range('D21').Value = '''6.73'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = '''78.51'
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('D23').Value = '''254.89'
$ws.Range('E23').Value = '  -1.53%  '
$ws.Range('D24').Value = '''1.94'
$ws.Range('E24').Value = '  +8.76%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').Value = '''2.50'
$ws.Range('E27').Value = '  -3.40%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '''10.51'
$ws.Range('E28').Value = '  -4.06%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.28'
$ws.Range('E29').Value = '  -1.28%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '''177.06'
$ws.Range('E30').Value = '  +1.06%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '''22.51'
$ws.Range('E31').Value = '  -3.49%  '
$ws.Range('E32').Value = '  -2.62%  '
$ws.Range('E33').Value = '  +0.77%  '
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('E35').Value = '  -4.63%  '
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range('E37').Value = '  -3.67%  '
$ws.Range('D38').Value = '''6.45'
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('D39').Value = '''2.39'
$ws.Range('E39').Value = '  -4.97%  '
$ws.Range('D40').Value = '''0.0276'
$ws.Range('E40').Value = '  -0.97%  '
$ws.Range('D41').Value = '''68.01'
$ws.Range('E41').Value = '  +24.55%  '
$ws.Range('D42').Value = '''5.14'
$ws.Range('E42').Value = '  +14.48%  '
$ws.Range('E43').Value = '  +8.22%  '
$ws.Range('D44').Value = '''9.21'
$ws.Range('E44').Value = '  +0.47%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '''18.89'
$ws.Range('E45').Value = '  -1.38%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '''0.201'
$ws.Range('E46').Value = '  +2.60%  '
$ws.Range('D47').Value = '''2.50'
$ws.Range('E47').Value = '  -1.91%  '
$ws.Range('E48').Value = '  -2.48%  '
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('D50').Value = '''99.84'
$ws.Range('E50').Value = '  -2.91%  '
$ws.Range('D51').Value = '''1.15'
$ws.Range('E51').Value = '  -5.73%  '
